$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 (Marion Gilbert) phone numbers L6:P6
$ws.Range("L6").Value = "(252) 232-3724"
$ws.Range("M6").Value = "(757) 424-0364"
$ws.Range("N6").Value = "(757) 424-1089"
$ws.Range("O6").Value = "(757) 482-3605"
$ws.Range("P6").Value = "(757) 560-6109"

# Row 7 (Mike Thompson) phone numbers L7:P7
$ws.Range("L7").Value = "(407) 344-9360"
$ws.Range("M7").Value = "(919) 402-4322"
$ws.Range("N7").Value = "(434) 791-2763"
$ws.Range("O7").Value = "(336) 602-2970"
$ws.Range("P7").Value = "(919) 926-0149"

# Row 10 (Heather Carter) phone numbers L10:P10
$ws.Range("L10").Value = "(321) 200-7555"
$ws.Range("M10").Value = "(904) 460-7673"
$ws.Range("N10").Value = "(904) 683-3096"
$ws.Range("O10").Value = "(904) 388-9866"
$ws.Range("P10").Value = "(904) 683-8680"
